# ---------------------------------------------------------------------------
# Applies the "import export and getPoints all work" commit to the workbook.
#
#   * Renames "Entries Base" -> "Entry Base"
#   * Renames "Tags"         -> "Tag Defs" (and adds _emoji/_desc columns)
#   * Adds a brand-new "Tags" worksheet at the end (did/pid -> tid mapping)
#   * Adds an _eid column to "Entry Points" (new point-entries schema)
#   * Adds _eid / _period columns to "Entry Base"
#   * Replaces/extends the sample data rows on "Defs" and "Point Defs"
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet renames + new sheet
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Entries Base").Name = "Entry Base"
$wb.Worksheets.Item("Tags").Name = "Tag Defs"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tagsSheet = $wb.Worksheets.Add($null, $lastSheet)
$tagsSheet.Name = "Tags"

# ---------------------------------------------------------------------------
# 2) "Defs" sheet - replace rows 2-5 with refreshed data, add rows 6-7
# ---------------------------------------------------------------------------
$defs = $wb.Worksheets.Item("Defs")

$defsRows = @(
    @("lgqy9rb3-0bmn", "2023-04-21T14:34:06.217", "lgqy9rbd", "'FALSE", "0m7w", "test one", "1️⃣", "Initial desc", "SECOND"),
    @("lgqy9rbd-avpc", "2023-04-21T14:34:06.217", "lgqy9rbe", "'FALSE", "ay7l", "twooo", "2️⃣", "now with a description", "WEEK"),
    @("lgqy9rbe-2ban", "2023-04-21T14:34:06.218", "lgr25h5o", "'TRUE", "05a8", "afree", "3️⃣", "Set a description", "SECOND"),
    @("lgqy9rbe-0keb", "2023-04-21T14:34:06.218", "lgqy9rbe", "'FALSE", "e0bq", "FOUR", "4️⃣", "having fun", "SECOND"),
    @("lgr0q1t4-3lil", "2023-04-21T15:42:45.553", "lgr0q1te", "'FALSE", "05a8", "afree", "3️⃣", "Edited with description!", "SECOND"),
    @("lgr0q1te-5odh", "2023-04-21T15:42:45.554", "lgr0q1te", "'FALSE", "7gor", "Five", "5️⃣", "not in first file, added to second", "SECOND")
)

$r = 2
foreach ($row in $defsRows) {
    $c = 1
    foreach ($val in $row) {
        $defs.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3) "Point Defs" sheet - add rows 2-4 (boolean column stays a real bool)
# ---------------------------------------------------------------------------
$pointDefs = $wb.Worksheets.Item("Point Defs")

$pointDefsRows = @(
    @("lgqy9rbe-0bcq", "2023-04-21T14:34:06.218", "lgr25h5o", $true,  "e0bq", "0pc6", "set alternatively", "☝️", "Set a description", "BOOL", "COUNT", "TEXT"),
    @("lgqy9rbe-3tnn", "2023-04-21T14:34:06.218", "lgqy9rbe", $false, "e0bq", "0tb7", "test point",         "🆕", "Set a description", "TEXT", "COUNT", "TEXT"),
    @("lgr0q1te-9rqg", "2023-04-21T15:42:45.555", "lgr0q1tf", $false, "e0bq", "0pc6", "updated label",      "☝️", "Set a description", "BOOL", "COUNT", "TEXT")
)

$r = 2
foreach ($row in $pointDefsRows) {
    $c = 1
    foreach ($val in $row) {
        $pointDefs.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 4) "Entry Base" sheet - insert _eid / _period columns before _note
# ---------------------------------------------------------------------------
$entryBase = $wb.Worksheets.Item("Entry Base")
$entryBase.Range("F1").Value = "_eid"
$entryBase.Range("G1").Value = "_period"
$entryBase.Range("H1").Value = "_note"

# ---------------------------------------------------------------------------
# 5) "Entry Points" sheet - insert _eid column before _val
# ---------------------------------------------------------------------------
$entryPoints = $wb.Worksheets.Item("Entry Points")
$entryPoints.Range("G1").Value = "_eid"
$entryPoints.Range("H1").Value = "_val"

# ---------------------------------------------------------------------------
# 6) "Tag Defs" sheet - append _emoji / _desc columns
# ---------------------------------------------------------------------------
$tagDefs = $wb.Worksheets.Item("Tag Defs")
$tagDefs.Range("G1").Value = "_emoji"
$tagDefs.Range("H1").Value = "_desc"

# ---------------------------------------------------------------------------
# 7) "Tags" sheet (new) - header row
# ---------------------------------------------------------------------------
$tags = $wb.Worksheets.Item("Tags")
$tagsHeaders = @("_uid", "_created", "_updated", "_deleted", "_did", "_pid", "tid")
$c = 1
foreach ($val in $tagsHeaders) {
    $tags.Cells.Item(1, $c).Value = $val
    $c = $c + 1
}
